$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------------
# 1) Title: "s Waterfall model vs. Incremental" -> "s V-Model vs. Agile"
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1).Range
$titleXml = '<w:p xmlns:w="' + $wNs + '">' +
  '<w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>SDLC&#8217;</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>s V-Model vs. Agile</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> model</w:t></w:r>' +
  '</w:p>'
$titlePara.InsertXML($titleXml)

# ---------------------------------------------------------------------------
# 2) Date line: "Date: 1/28/16" -> "Date: 3/19" + "/16" (two runs)
# ---------------------------------------------------------------------------
$datePara = $d.Paragraphs(7).Range
$dateXml = '<w:p xmlns:w="' + $wNs + '">' +
  '<w:pPr><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t>Date: 3/19</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t>/16</w:t></w:r>' +
  '</w:p>'
$datePara.InsertXML($dateXml)

# ---------------------------------------------------------------------------
# 3) Big SDLC-definition paragraph (paragraph 10): Waterfall/Incremental -> Agile/V-model
# ---------------------------------------------------------------------------
$defPara = $d.Paragraphs(10).Range
$defXml = '<w:p xmlns:w="' + $wNs + '">' +
  '<w:pPr><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t>SDLC (</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t>Software Development Life Cycle) is a term used to describe the process for planning, creating, testing, and deploying an information system. This concept applies to a range of hardware and software configurations, as a system can be composed of hardware only, software only, or a combination of the two. There are two different models that I will discuss in th</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t>is paper one being the agile</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t xml:space="preserve"> model and </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t>the other being the V-model. The Agile model shows rapid and continuous improvement through the collaboration of small self-organized teams.</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t xml:space="preserve"> The </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t>V-</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t xml:space="preserve">model </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t>demonstrates relationships between each phase of the development life cycle and its associated phase of testing</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t>.</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '</w:p>'
$defPara.InsertXML($defXml)

# ---------------------------------------------------------------------------
# 4) Waterfall paragraph (12) -> Agile description paragraph
# ---------------------------------------------------------------------------
$agilePara = $d.Paragraphs(12).Range
$agileXml = '<w:p xmlns:w="' + $wNs + '">' +
  '<w:pPr><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:tab/><w:t>T</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t>he Agile software development process is a set of principles in which</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t xml:space="preserve"> requirements and solutions are furthered through </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t>collaboration betwe</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t>en small self-organizing, cross-</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t xml:space="preserve">functional teams. </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t xml:space="preserve">It promotes adaptive planning, outstanding development, quick deliveries, endless improvements, and it encourages quick and flexible response to change. </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t xml:space="preserve">Though Agile does not have specific methods to achieve these goals, it is evident that many have in fact grown a great deal from this and are now considered </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light"/><w:i/></w:rPr><w:t>Agile.</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t xml:space="preserve">The agile manifesto is based on twelve principles that mainly revolve around hard work, self-determination, close </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t>face-to-face</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t xml:space="preserve"> interactions, and </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t>simplicity.</w:t></w:r>' +
  '</w:p>'
$agilePara.InsertXML($agileXml)

# ---------------------------------------------------------------------------
# 5) Incremental paragraph (14) -> V-model description paragraph (keep bookmark)
# ---------------------------------------------------------------------------
$vmodelPara = $d.Paragraphs(14).Range
$vmodelXml = '<w:p xmlns:w="' + $wNs + '">' +
  '<w:pPr><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:tab/><w:t xml:space="preserve">The </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t>V-model represents a development process that can be considered an extension of the waterfall model, and is an example of the less intricate V-model.</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t xml:space="preserve"> Instead of going straight down in a linear way, the process steps are bent after its coding phase to form a shape that resembles a V.</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/></w:rPr><w:t>The horizontal and vertical axes represent the time or project completeness and level of abstraction, which is read from left to right. This model has 4 validation phases, that being Unit testing, Integration testing, System testing, and User acceptance testing.</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '</w:p>'
$vmodelPara.InsertXML($vmodelXml)
